# The "Förändrad" (Changed) column C holds a date serial number that was
# bumped by one day (2026-02-28 -> 2026-03-01, i.e. serial 46081 -> 46082)
# for every data row (rows 2-287) in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C287").Value = 46082
